$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3001.5
$ws.Range("I86").Value = 3003
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 3003
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1880
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 3001.5
$ws.Range("I89").Value = 3003
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 15015
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -9399
$ws.Range("N89").Value = -26232
$ws.Range("H106").Value = 22423402
$ws.Range("I106").Value = 22423402
$ws.Range("K106").Value = 22423402
$ws.Range("M106").Value = -22422771
$ws.Range("H134").Value = 58000
$ws.Range("J134").Value = 58000
$ws.Range("L134").Value = 58000
$ws.Range("N134").Value = -68140
$ws.Range("H138").Value = 5683450.5
$ws.Range("I138").Value = 701.44116
$ws.Range("J138").Value = 9261478
$ws.Range("K138").Value = 2104.32348
$ws.Range("L138").Value = 27784434
$ws.Range("M138").Value = 3035.67652
$ws.Range("N138").Value = -27794714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2107
$ws.Range("I45").Value = 1408
$ws.Range("K45").Value = 1408
$ws.Range("M45").Value = -1031
$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50676
$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52340
$ws.Range("H97").Value = 16667206
$ws.Range("I97").Value = 30303500
$ws.Range("J97").Value = 624.6667
$ws.Range("K97").Value = 30303500
$ws.Range("L97").Value = 624.6667
$ws.Range("M97").Value = -30303004
$ws.Range("N97").Value = -1616.6667
$ws.Range("H122").Value = 6844.316
$ws.Range("I122").Value = 9100.77
$ws.Range("J122").Value = 1955.3334
$ws.Range("K122").Value = 27302.31
$ws.Range("L122").Value = 5866.0002
$ws.Range("M122").Value = -24852.31
$ws.Range("N122").Value = -10766.0002
$ws.Range("H132").Value = 2031.3529
$ws.Range("I132").Value = 1467.4
$ws.Range("J132").Value = 4082.0908
$ws.Range("K132").Value = 4402.200000000001
$ws.Range("L132").Value = 12246.2724
$ws.Range("M132").Value = -1872.200000000001
$ws.Range("N132").Value = -17306.2724
$ws.Range("H139").Value = 66614.164
$ws.Range("J139").Value = 66614.164
$ws.Range("L139").Value = 66614.164
$ws.Range("N139").Value = -76894.164

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1263.2069
$ws.Range("I20").Value = 917.2941
$ws.Range("J20").Value = 1753.25
$ws.Range("K20").Value = 917.2941
$ws.Range("L20").Value = 1753.25
$ws.Range("M20").Value = -670.2941
$ws.Range("N20").Value = -2247.25
$ws.Range("H107").Value = 1433.2354
$ws.Range("I107").Value = 1166.5454
$ws.Range("J107").Value = 1922.1666
$ws.Range("K107").Value = 1166.5454
$ws.Range("L107").Value = 1922.1666
$ws.Range("M107").Value = 753.4546
$ws.Range("N107").Value = -5762.1666
$ws.Range("H118").Value = 27963
$ws.Range("J118").Value = 27963
$ws.Range("L118").Value = 27963
$ws.Range("N118").Value = -31277
$ws.Range("H132").Value = 12225
$ws.Range("J132").Value = 12225
$ws.Range("L132").Value = 12225
$ws.Range("N132").Value = -22345

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 78.625
$ws.Range("I7").Value = 61.285713
$ws.Range("K7").Value = 61.285713
$ws.Range("M7").Value = 51.714287
$ws.Range("I17").Value = 50000
$ws.Range("J17").Value = 5000
$ws.Range("K17").Value = 50000
$ws.Range("L17").Value = 5000
$ws.Range("M17").Value = -49826
$ws.Range("N17").Value = -5348
$ws.Range("H22").Value = 248.88889
$ws.Range("I22").Value = 252.42857
$ws.Range("J22").Value = 236.5
$ws.Range("K22").Value = 252.42857
$ws.Range("L22").Value = 236.5
$ws.Range("M22").Value = 97.57142999999999
$ws.Range("N22").Value = -936.5
$ws.Range("H31").Value = 1199.15
$ws.Range("I31").Value = 902.36957
$ws.Range("J31").Value = 1451.963
$ws.Range("K31").Value = 902.36957
$ws.Range("L31").Value = 1451.963
$ws.Range("M31").Value = -607.36957
$ws.Range("N31").Value = -2041.963
$ws.Range("H34").Value = 1199.15
$ws.Range("I34").Value = 902.36957
$ws.Range("J34").Value = 1451.963
$ws.Range("K34").Value = 902.36957
$ws.Range("L34").Value = 1451.963
$ws.Range("M34").Value = -700.36957
$ws.Range("N34").Value = -1855.963
$ws.Range("H132").Value = 2019.0454
$ws.Range("I132").Value = 1427.2703
$ws.Range("J132").Value = 5147
$ws.Range("K132").Value = 4281.810899999999
$ws.Range("L132").Value = 15441
$ws.Range("M132").Value = -1751.810899999999
$ws.Range("N132").Value = -20501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 411.75
$ws.Range("I11").Value = 411.75
$ws.Range("K11").Value = 1235.25
$ws.Range("M11").Value = -1095.25
$ws.Range("H39").Value = 8360.757
$ws.Range("J39").Value = 8360.757
$ws.Range("L39").Value = 25082.271
$ws.Range("N39").Value = -25670.271
$ws.Range("H68").Value = 880.55
$ws.Range("I68").Value = 738.7538500000001
$ws.Range("J68").Value = 1143.8857
$ws.Range("K68").Value = 2216.26155
$ws.Range("L68").Value = 3431.6571
$ws.Range("M68").Value = -1405.26155
$ws.Range("N68").Value = -5053.6571
$ws.Range("H71").Value = 880.55
$ws.Range("I71").Value = 738.7538500000001
$ws.Range("J71").Value = 1143.8857
$ws.Range("K71").Value = 6648.784650000001
$ws.Range("L71").Value = 10294.9713
$ws.Range("M71").Value = -2592.784650000001
$ws.Range("N71").Value = -18406.9713
$ws.Range("H80").Value = 1200
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1200
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3600
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -5472
$ws.Range("H83").Value = 1200
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1200
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 10800
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -20160
$ws.Range("H92").Value = 800
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H107").Value = 935.42
$ws.Range("I107").Value = 279.9091
$ws.Range("J107").Value = 1120.3077
$ws.Range("K107").Value = 839.7273
$ws.Range("L107").Value = 3360.9231
$ws.Range("M107").Value = 1080.2727
$ws.Range("N107").Value = -7200.9231
$ws.Range("H113").Value = 15152007
$ws.Range("I113").Value = 481
$ws.Range("J113").Value = 22727770
$ws.Range("K113").Value = 1443
$ws.Range("L113").Value = 68183310
$ws.Range("M113").Value = 727
$ws.Range("N113").Value = -68187650
$ws.Range("H122").Value = 723.5625
$ws.Range("I122").Value = 630.4666999999999
$ws.Range("J122").Value = 805.7059
$ws.Range("K122").Value = 5674.2003
$ws.Range("L122").Value = 7251.3531
$ws.Range("M122").Value = -3224.2003
$ws.Range("N122").Value = -12151.3531
$ws.Range("H123").Value = 1500
$ws.Range("I123").Value = 1500
$ws.Range("K123").Value = 4500
$ws.Range("M123").Value = -2050
$ws.Range("H132").Value = 1164.0667
$ws.Range("I132").Value = 1025.3334
$ws.Range("J132").Value = 1198.75
$ws.Range("K132").Value = 9228.000599999999
$ws.Range("L132").Value = 10788.75
$ws.Range("M132").Value = -6698.000599999999
$ws.Range("N132").Value = -15848.75
$ws.Range("H134").Value = 15181.758
$ws.Range("I134").Value = 1819.8
$ws.Range("J134").Value = 20991.305
$ws.Range("K134").Value = 5459.4
$ws.Range("L134").Value = 62973.915
$ws.Range("M134").Value = -389.3999999999996
$ws.Range("N134").Value = -73113.91500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 775.7353000000001
$ws.Range("I97").Value = 677.64703
$ws.Range("J97").Value = 873.82355
$ws.Range("K97").Value = 677.64703
$ws.Range("L97").Value = 873.82355
$ws.Range("M97").Value = -181.64703
$ws.Range("N97").Value = -1865.82355
$ws.Range("H102").Value = 1875.8334
$ws.Range("I102").Value = 1269.2084
$ws.Range("J102").Value = 4302.3335
$ws.Range("K102").Value = 1269.2084
$ws.Range("L102").Value = 4302.3335
$ws.Range("M102").Value = 352.7916
$ws.Range("N102").Value = -7546.3335
$ws.Range("H122").Value = 618358.3
$ws.Range("I122").Value = 1111906.4
$ws.Range("J122").Value = 1423.25
$ws.Range("K122").Value = 3335719.2
$ws.Range("L122").Value = 4269.75
$ws.Range("M122").Value = -3333269.2
$ws.Range("N122").Value = -9169.75
$ws.Range("H126").Value = 2049.0715
$ws.Range("I126").Value = 1536.4615
$ws.Range("J126").Value = 2493.3333
$ws.Range("K126").Value = 4609.3845
$ws.Range("L126").Value = 7479.999899999999
$ws.Range("M126").Value = -2139.3845
$ws.Range("N126").Value = -12419.9999
$ws.Range("H132").Value = 2597.302
$ws.Range("I132").Value = 2127.561
$ws.Range("J132").Value = 4202.25
$ws.Range("K132").Value = 6382.683000000001
$ws.Range("L132").Value = 12606.75
$ws.Range("M132").Value = -3852.683000000001
$ws.Range("N132").Value = -17666.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4071.7188
$ws.Range("I132").Value = 3392.6875
$ws.Range("J132").Value = 4750.75
$ws.Range("K132").Value = 10178.0625
$ws.Range("L132").Value = 14252.25
$ws.Range("M132").Value = -7648.0625
$ws.Range("N132").Value = -19312.25
$ws.Range("H136").Value = 5450.1665
$ws.Range("I136").Value = 2500.1667
$ws.Range("J136").Value = 8400.166999999999
$ws.Range("K136").Value = 7500.500100000001
$ws.Range("L136").Value = 25200.501
$ws.Range("M136").Value = -4950.500100000001
$ws.Range("N136").Value = -30300.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 40995.32
$ws.Range("I122").Value = 53488.58
$ws.Range("J122").Value = 1433.3334
$ws.Range("K122").Value = 160465.74
$ws.Range("L122").Value = 4300.0002
$ws.Range("M122").Value = -158015.74
$ws.Range("N122").Value = -9200.0002
$ws.Range("H132").Value = 16131341
$ws.Range("I132").Value = 22729178
$ws.Range("J132").Value = 3293.2222
$ws.Range("K132").Value = 68187534
$ws.Range("L132").Value = 9879.6666
$ws.Range("M132").Value = -68185004
$ws.Range("N132").Value = -14939.6666
$ws.Range("H136").Value = 11941781
$ws.Range("I136").Value = 16717296
$ws.Range("J136").Value = 2993.875
$ws.Range("K136").Value = 50151888
$ws.Range("L136").Value = 8981.625
$ws.Range("M136").Value = -50149338
$ws.Range("N136").Value = -14081.625
